$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename activity description in C40 from "Theme Party" to "Kolokium Zon Selatan"
$ws.Range("C40").Value = "Kolokium Zon Selatan"

# 2. Update the merit value for that activity (D40) from 100 to 0
$ws.Range("D40").Value = 0

# 3. Merge B15:C15 (to match B13:C13 / B14:C14 above it)
$ws.Range("B15:C15").Merge()
